# Auto-generated script applying numeric updates described by the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (30 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 20027870
$ws.Range("I40").Value = 16198.857
$ws.Range("J40").Value = 45497270
$ws.Range("K40").Value = 16198.857
$ws.Range("L40").Value = 45497270
$ws.Range("M40").Value = -16023.857
$ws.Range("N40").Value = -45497620
$ws.Range("H87").Value = 85000
$ws.Range("J87").Value = 85000
$ws.Range("L87").Value = 85000
$ws.Range("N87").Value = -87496
$ws.Range("H90").Value = 85000
$ws.Range("J90").Value = 85000
$ws.Range("L90").Value = 255000
$ws.Range("N90").Value = -267480
$ws.Range("H100").Value = 1556.5
$ws.Range("I100").Value = 1556.5
$ws.Range("K100").Value = 1556.5
$ws.Range("M100").Value = -1015.5
$ws.Range("H113").Value = 2969.423
$ws.Range("I113").Value = 1948.7778
$ws.Range("J113").Value = 3509.7646
$ws.Range("K113").Value = 1948.7778
$ws.Range("L113").Value = 3509.7646
$ws.Range("M113").Value = 1305.2222
$ws.Range("N113").Value = -10017.7646
$ws.Range("H132").Value = 138126.5
$ws.Range("I132").Value = 252781.22
$ws.Range("K132").Value = 758343.66
$ws.Range("M132").Value = -755813.66

# ---- Sheet: ARM (24 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3747.182
$ws.Range("I32").Value = 2032.2153
$ws.Range("K32").Value = 2032.2153
$ws.Range("M32").Value = -1745.2153
$ws.Range("H37").Value = 49999
$ws.Range("J37").Value = 49999
$ws.Range("L37").Value = 49999
$ws.Range("N37").Value = -50545
$ws.Range("H52").Value = 69999
$ws.Range("J52").Value = 69999
$ws.Range("L52").Value = 69999
$ws.Range("N52").Value = -70525
$ws.Range("H61").Value = 31326.8
$ws.Range("I61").Value = 31326.8
$ws.Range("K61").Value = 31326.8
$ws.Range("M61").Value = -31114.8
$ws.Range("H102").Value = 572618.4399999999
$ws.Range("I102").Value = 762559.1
$ws.Range("K102").Value = 762559.1
$ws.Range("M102").Value = -760937.1
$ws.Range("H136").Value = 31326.8
$ws.Range("I136").Value = 31326.8
$ws.Range("K136").Value = 93980.39999999999
$ws.Range("M136").Value = -91430.39999999999

# ---- Sheet: BSM (20 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 24000
$ws.Range("I88").Value = 24000
$ws.Range("K88").Value = 24000
$ws.Range("M88").Value = -23594
$ws.Range("H91").Value = 24000
$ws.Range("I91").Value = 24000
$ws.Range("K91").Value = 24000
$ws.Range("M91").Value = -22596
$ws.Range("H103").Value = 64828.5
$ws.Range("J103").Value = 64828.5
$ws.Range("L103").Value = 64828.5
$ws.Range("N103").Value = -67172.5
$ws.Range("H121").Value = 69999
$ws.Range("J121").Value = 69999
$ws.Range("L121").Value = 69999
$ws.Range("N121").Value = -73493
$ws.Range("H133").Value = 112063.25
$ws.Range("J133").Value = 112063.25
$ws.Range("L133").Value = 112063.25
$ws.Range("N133").Value = -122183.25

# ---- Sheet: CRP (15 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5068.857
$ws.Range("I94").Value = 4450
$ws.Range("J94").Value = 5316.4
$ws.Range("K94").Value = 4450
$ws.Range("L94").Value = 5316.4
$ws.Range("M94").Value = -3999
$ws.Range("N94").Value = -6218.4
$ws.Range("H124").Value = 31442
$ws.Range("J124").Value = 31442
$ws.Range("L124").Value = 31442
$ws.Range("N124").Value = -36352
$ws.Range("H131").Value = 48848.5
$ws.Range("J131").Value = 48848.5
$ws.Range("L131").Value = 48848.5
$ws.Range("N131").Value = -58928.5

# ---- Sheet: CUL (34 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 278.08
$ws.Range("J2").Value = 344.58334
$ws.Range("L2").Value = 2067.50004
$ws.Range("N2").Value = -2293.50004
$ws.Range("H8").Value = 83.333336
$ws.Range("I8").Value = 83.333336
$ws.Range("K8").Value = 250.000008
$ws.Range("M8").Value = -111.000008
$ws.Range("H9").Value = 918.4
$ws.Range("I9").Value = 899
$ws.Range("K9").Value = 2697
$ws.Range("M9").Value = -2473
$ws.Range("H38").Value = 205.72223
$ws.Range("J38").Value = 248.07692
$ws.Range("L38").Value = 744.23076
$ws.Range("N38").Value = -1438.23076
$ws.Range("H74").Value = 14645.111
$ws.Range("I74").Value = 3887.75
$ws.Range("J74").Value = 23251
$ws.Range("K74").Value = 11663.25
$ws.Range("L74").Value = 69753
$ws.Range("M74").Value = -10602.25
$ws.Range("N74").Value = -71875
$ws.Range("H77").Value = 14645.111
$ws.Range("I77").Value = 3887.75
$ws.Range("J77").Value = 23251
$ws.Range("K77").Value = 34989.75
$ws.Range("L77").Value = 209259
$ws.Range("M77").Value = -29685.75
$ws.Range("N77").Value = -219867
$ws.Range("H122").Value = 788.5
$ws.Range("I122").Value = 788.5
$ws.Range("K122").Value = 7096.5
$ws.Range("M122").Value = -4646.5

# ---- Sheet: GSM (31 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 27499.25
$ws.Range("J45").Value = 27499.25
$ws.Range("L45").Value = 27499.25
$ws.Range("N45").Value = -28617.25
$ws.Range("H80").Value = 93227.28999999999
$ws.Range("I80").Value = 256399.5
$ws.Range("K80").Value = 256399.5
$ws.Range("M80").Value = -255401.5
$ws.Range("H83").Value = 93227.28999999999
$ws.Range("I83").Value = 256399.5
$ws.Range("K83").Value = 1281997.5
$ws.Range("M83").Value = -1277005.5
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H107").Value = 7937340
$ws.Range("I107").Value = 23810024
$ws.Range("J107").Value = 997.5
$ws.Range("K107").Value = 23810024
$ws.Range("L107").Value = 997.5
$ws.Range("M107").Value = -23808104
$ws.Range("N107").Value = -4837.5
$ws.Range("H127").Value = 80250.5
$ws.Range("J127").Value = 80250.5
$ws.Range("L127").Value = 80250.5
$ws.Range("N127").Value = -90170.5
$ws.Range("H134").Value = 34809.832
$ws.Range("J134").Value = 34809.832
$ws.Range("L134").Value = 104429.496
$ws.Range("N134").Value = -109499.496

# ---- Sheet: LTW (53 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4275.636
$ws.Range("I7").Value = 2890.2
$ws.Range("K7").Value = 2890.2
$ws.Range("M7").Value = -2778.2
$ws.Range("H16").Value = 6453068.5
$ws.Range("I16").Value = 9525123
$ws.Range("J16").Value = 1753.4
$ws.Range("K16").Value = 9525123
$ws.Range("L16").Value = 1753.4
$ws.Range("M16").Value = -9524953
$ws.Range("N16").Value = -2093.4
$ws.Range("H22").Value = 1164.6471
$ws.Range("I22").Value = 1011.1111
$ws.Range("K22").Value = 1011.1111
$ws.Range("M22").Value = -716.1111
$ws.Range("H27").Value = 1164.6471
$ws.Range("I27").Value = 1011.1111
$ws.Range("K27").Value = 1011.1111
$ws.Range("M27").Value = -904.1111
$ws.Range("H46").Value = 6064.9062
$ws.Range("I46").Value = 3666.6667
$ws.Range("J46").Value = 6313
$ws.Range("K46").Value = 3666.6667
$ws.Range("L46").Value = 6313
$ws.Range("M46").Value = -3478.6667
$ws.Range("N46").Value = -6689
$ws.Range("H55").Value = 273.41666
$ws.Range("I55").Value = 122.666664
$ws.Range("K55").Value = 122.666664
$ws.Range("M55").Value = 50.333336
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H93").Value = 1650.4
$ws.Range("I93").Value = 1375.5
$ws.Range("J93").Value = 2750
$ws.Range("K93").Value = 1375.5
$ws.Range("L93").Value = 2750
$ws.Range("M93").Value = -127.5
$ws.Range("N93").Value = -5246
$ws.Range("H126").Value = 4275.636
$ws.Range("I126").Value = 2890.2
$ws.Range("K126").Value = 8670.599999999999
$ws.Range("M126").Value = -6200.599999999999
$ws.Range("H136").Value = 5274.0625
$ws.Range("I136").Value = 2988.5
$ws.Range("K136").Value = 8965.5
$ws.Range("M136").Value = -6415.5

# ---- Sheet: WVR (38 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H62").Value = 9776.333000000001
$ws.Range("I62").Value = 4664.5
$ws.Range("K62").Value = 4664.5
$ws.Range("M62").Value = -4040.5
$ws.Range("H65").Value = 9776.333000000001
$ws.Range("I65").Value = 4664.5
$ws.Range("K65").Value = 23322.5
$ws.Range("M65").Value = -20202.5
$ws.Range("H107").Value = 2352.7812
$ws.Range("I107").Value = 2289.1738
$ws.Range("K107").Value = 6867.5214
$ws.Range("M107").Value = -4947.5214
$ws.Range("H122").Value = 2711.152
$ws.Range("I122").Value = 2863.158
$ws.Range("K122").Value = 8589.474
$ws.Range("M122").Value = -6139.474
$ws.Range("H126").Value = 2878.4375
$ws.Range("I126").Value = 2492.1667
$ws.Range("J126").Value = 4037.25
$ws.Range("K126").Value = 7476.500100000001
$ws.Range("L126").Value = 12111.75
$ws.Range("M126").Value = -5006.500100000001
$ws.Range("N126").Value = -17051.75
$ws.Range("H132").Value = 61729450
$ws.Range("I132").Value = 7937149
$ws.Range("J132").Value = 250002500
$ws.Range("K132").Value = 23811447
$ws.Range("L132").Value = 750007500
$ws.Range("M132").Value = -23808917
$ws.Range("N132").Value = -750012560
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

